$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on Price cells whose new value would otherwise be
# auto-parsed as a number (losing formatting such as trailing zeros), so the
# written value stays a plain text string like the original inline strings.
$forceTextCells = @("D4","D5","D6","D7","D8","D9","D10","D11","D14","D15","D16","D18","D19","D20","D21","D22","D24","D26","D27","D28","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D47","D48","D49","D50","D51")
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Cell value updates (Price + Volume(1h) columns), row by row.
$ws.Range("D2").Value = "20.401.41"
$ws.Range("E2").Value = "  +2.30%  "
$ws.Range("D3").Value = "1.460.89"
$ws.Range("E3").Value = "  +3.35%  "
$ws.Range("D4").Value = "1.009"
$ws.Range("E4").Value = "  +0.83%  "
$ws.Range("D5").Value = "0.9501"
$ws.Range("D6").Value = "274.76"
$ws.Range("E6").Value = "  -0.42%  "
$ws.Range("D7").Value = "0.3654"
$ws.Range("E7").Value = "  -0.27%  "
$ws.Range("D8").Value = "0.3064"
$ws.Range("E8").Value = "  -0.94%  "
$ws.Range("D9").Value = "39.77"
$ws.Range("E9").Value = "  +0.17%  "
$ws.Range("D10").Value = "1.034"
$ws.Range("E10").Value = "  +0.29%  "
$ws.Range("D11").Value = "0.06574"
$ws.Range("E11").Value = "  +0.88%  "
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("E13").Value = "  -1.16%  "
$ws.Range("D14").Value = "17.89"
$ws.Range("E14").Value = "  +2.04%  "
$ws.Range("D15").Value = "6.135"
$ws.Range("E15").Value = "  -0.75%  "
$ws.Range("D16").Value = "0.00001022"
$ws.Range("E16").Value = "  +0.62%  "
$ws.Range("D17").Value = "1.458.89"
$ws.Range("E17").Value = "  +3.15%  "
$ws.Range("D18").Value = "0.9677"
$ws.Range("E18").Value = "  -3.29%  "
$ws.Range("D19").Value = "0.05803"
$ws.Range("E19").Value = "  +2.56%  "
$ws.Range("D20").Value = "69.40"
$ws.Range("E20").Value = "  -2.34%  "
$ws.Range("D21").Value = "5.431"
$ws.Range("E21").Value = "  -3.13%  "
$ws.Range("D22").Value = "14.45"
$ws.Range("E22").Value = "  -1.52%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").Value = "2.246"
$ws.Range("E24").Value = "  +0.64%  "
$ws.Range("D25").Value = "20.427.44"
$ws.Range("E25").Value = "  +2.36%  "
$ws.Range("D26").Value = "141.42"
$ws.Range("E26").Value = "  +6.44%  "
$ws.Range("D27").Value = "2.080"
$ws.Range("E27").Value = "  -7.65%  "
$ws.Range("D28").Value = "17.11"
$ws.Range("E28").Value = "  -0.63%  "
$ws.Range("D29").Value = "1.612.94"
$ws.Range("E29").Value = "  +2.65%  "
$ws.Range("D30").Value = "112.21"
$ws.Range("E30").Value = "  +2.48%  "
$ws.Range("D31").Value = "3.817"
$ws.Range("E31").Value = "  -1.85%  "
$ws.Range("D32").Value = "4.887"
$ws.Range("E32").Value = "  -7.15%  "
$ws.Range("D33").Value = "0.07888"
$ws.Range("E33").Value = "  +2.69%  "
$ws.Range("D34").Value = "0.7897"
$ws.Range("E34").Value = "  -3.28%  "
$ws.Range("D35").Value = "1.532"
$ws.Range("E35").Value = "  +2.86%  "
$ws.Range("D36").Value = "0.05723"
$ws.Range("E36").Value = "  -0.84%  "
$ws.Range("D37").Value = "1.152"
$ws.Range("E37").Value = "  +5.34%  "
$ws.Range("D38").Value = "4.678"
$ws.Range("E38").Value = "  -4.60%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.02022"
$ws.Range("E39").Value = "  -1.27%  "
$ws.Range("B40").Value = "Frax"
$ws.Range("C40").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D40").Value = "0.9563"
$ws.Range("E40").Value = "  -4.27%  "
$ws.Range("D41").Value = "10.32"
$ws.Range("E41").Value = "  -0.87%  "
$ws.Range("D42").Value = "7.474"
$ws.Range("E42").Value = "  -10.00%  "
$ws.Range("D43").Value = "0.1856"
$ws.Range("E43").Value = "  -1.40%  "
$ws.Range("D44").Value = "0.5255"
$ws.Range("E44").Value = "  -0.79%  "
$ws.Range("D45").Value = "3.489"
$ws.Range("E45").Value = "  -1.31%  "
$ws.Range("E46").Value = "  -3.72%  "
$ws.Range("D47").Value = "116.93"
$ws.Range("E47").Value = "  +1.52%  "
$ws.Range("D48").Value = "0.5126"
$ws.Range("E48").Value = "  -0.88%  "
$ws.Range("D49").Value = "1.746"
$ws.Range("E49").Value = "  -1.06%  "
$ws.Range("D50").Value = "0.06413"
$ws.Range("E50").Value = "  +3.73%  "
$ws.Range("D51").Value = "0.9927"
$ws.Range("E51").Value = "  -0.82%  "
